$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 2.05
$ws.Range("O2").Value = 1.75

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3
$ws.Range("N3").Value = 2.63
$ws.Range("O3").Value = 1.5
$ws.Range("V3").Value = 10
$ws.Range("W3").Value = 23
$ws.Range("X3").Value = 23

# Row 4
$ws.Range("N4").Value = 2.1
$ws.Range("O4").Value = 1.73

# Row 5
$ws.Range("L5").Value = 1.25
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 1.83
$ws.Range("O5").Value = 2.03

# Row 6
$ws.Range("G6").Value = 1.37
$ws.Range("H6").Value = 4.45
$ws.Range("V6").Value = 7.1
$ws.Range("W6").Value = 8.75
$ws.Range("X6").Value = 8.75
$ws.Range("Z6").Value = 18
$ws.Range("AA6").Value = 8.5
$ws.Range("AB6").Value = 12.5
$ws.Range("AC6").Value = 37
$ws.Range("AD6").Value = 175

# Row 7
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 8
$ws.Range("R7").Value = 1.81
$ws.Range("S7").Value = 1.89
$ws.Range("T7").Value = 7.3
$ws.Range("V7").Value = 7.2
$ws.Range("Y7").Value = 18
$ws.Range("Z7").Value = 15.5
$ws.Range("AB7").Value = 15.5
$ws.Range("AC7").Value = 55
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 50
$ws.Range("AH7").Value = 150

# Row 8
$ws.Range("I8").Value = 3.3
$ws.Range("M8").Value = 4.25
$ws.Range("N8").Value = 1.65
$ws.Range("O8").Value = 2.2
$ws.Range("P8").Value = 1.32
$ws.Range("Q8").Value = 3.2
$ws.Range("R8").Value = 1.55
$ws.Range("S8").Value = 2.32
$ws.Range("U8").Value = 12
$ws.Range("AC8").Value = 34
$ws.Range("AE8").Value = 15
$ws.Range("AF8").Value = 21

# Row 11
$ws.Range("G11").Value = 1.75
$ws.Range("I11").Value = 4
$ws.Range("U11").Value = 9.5

# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 3.2
$ws.Range("N12").Value = 2.03
$ws.Range("O12").Value = 1.83
$ws.Range("P12").Value = 1.38
$ws.Range("Q12").Value = 2.87
$ws.Range("S12").Value = 2
$ws.Range("U12").Value = 11
$ws.Range("X12").Value = 19
$ws.Range("AA12").Value = 6.5
$ws.Range("AB12").Value = 15
$ws.Range("AD12").Value = 201

# Row 13
$ws.Range("I13").Value = 3.05
$ws.Range("N13").Value = 2.12
$ws.Range("O13").Value = 1.65
$ws.Range("U13").Value = 10.25
$ws.Range("AB13").Value = 15.5
$ws.Range("AC13").Value = 80

# Row 14
$ws.Range("G14").Value = 1.62
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 5
$ws.Range("L14").Value = 1.22
$ws.Range("M14").Value = 4
$ws.Range("N14").Value = 1.7
$ws.Range("O14").Value = 2.1
$ws.Range("U14").Value = 8.5
$ws.Range("X14").Value = 13
$ws.Range("AA14").Value = 7.5
$ws.Range("AE14").Value = 15

# Row 15
$ws.Range("G15").Value = 1.88
$ws.Range("H15").Value = 3.35
$ws.Range("I15").Value = 3.6
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 1.8
$ws.Range("T15").Value = 6.3
$ws.Range("U15").Value = 7.7
$ws.Range("V15").Value = 7.1
$ws.Range("W15").Value = 13.5
$ws.Range("X15").Value = 12
$ws.Range("Y15").Value = 20
$ws.Range("Z15").Value = 10
$ws.Range("AA15").Value = 5.8
$ws.Range("AB15").Value = 11.25
$ws.Range("AE15").Value = 9.5
$ws.Range("AF15").Value = 16.5
$ws.Range("AG15").Value = 10.25
$ws.Range("AH15").Value = 40
$ws.Range("AI15").Value = 25
$ws.Range("AJ15").Value = 28

# Row 16
$ws.Range("R16").Value = 2.21
$ws.Range("S16").Value = 1.6

# Row 17
$ws.Range("H17").Value = 3.1
$ws.Range("J17").Value = 1.08
$ws.Range("K17").Value = 8
$ws.Range("N17").Value = 2.25
$ws.Range("O17").Value = 1.62
$ws.Range("P17").Value = 1.5
$ws.Range("Q17").Value = 2.5
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("U17").Value = 10
$ws.Range("Z17").Value = 8
$ws.Range("AD17").Value = 351
$ws.Range("AE17").Value = 8.5

# Row 18
$ws.Range("P18").Value = 1.25

# Row 19
$ws.Range("R19").Value = 1.7
$ws.Range("S19").Value = 2.05

# Row 20
$ws.Range("P20").Value = 1.4

# Row 22
$ws.Range("T22").Value = 6.4
$ws.Range("U22").Value = 7.1
$ws.Range("V22").Value = 8
$ws.Range("AD22").Value = 800
$ws.Range("AE22").Value = 13.5

# Row 23
$ws.Range("G23").Value = 2.15
$ws.Range("H23").Value = 3.3
$ws.Range("J23").Value = 1.06
$ws.Range("K23").Value = 9.5
$ws.Range("L23").Value = 1.3
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 2.05
$ws.Range("O23").Value = 1.75
$ws.Range("R23").Value = 1.8
$ws.Range("S23").Value = 1.91
$ws.Range("T23").Value = 7.5
$ws.Range("V23").Value = 9
$ws.Range("W23").Value = 19
$ws.Range("X23").Value = 17
$ws.Range("Z23").Value = 9.5
$ws.Range("AA23").Value = 6.5
$ws.Range("AD23").Value = 251
$ws.Range("AE23").Value = 10
$ws.Range("AH23").Value = 41
$ws.Range("AJ23").Value = 34

